# The "dataset_meta_data" sheet had a redundant / duplicate column F
# ("official_cruise_name(s)") whose sole purpose was superseded by the
# existing "cruise_names" column further to the right. The commit removes
# that obsolete column entirely (shifting everything after it one column
# to the left, which also drops the now-unused shared strings for its
# header & its placeholder/description text), and then leaves the
# "dataset_meta_data" sheet as the active/selected sheet & tab (it was
# "vars_meta_data" before).

$wb = $excel.ActiveWorkbook

$dataset = $wb.Worksheets.Item("dataset_meta_data")

# Remove column F (official_cruise_name(s)) entirely; everything to the
# right (dataset_source ... cruise_names) shifts left by one column.
$dataset.Columns("F").Delete()

# Make "dataset_meta_data" the active sheet/tab, with B22 selected
# (previously "vars_meta_data" was active, with B28 selected there).
$dataset.Activate()
$dataset.Range("B22").Select()
